# Splits the single long <w:t> run in three paragraphs ("Programa resumido",
# "Programa" and "Bibliografia") into several <w:t> segments interleaved with
# manual line breaks (<w:br/>), matching the target OOXML diff exactly
# (including which segments keep xml:space="preserve").
#
# Find/Replace + plain text assignment in this runtime do not reproduce the
# per-segment xml:space="preserve" placement that the diff expects, so the
# paragraph's Range content is replaced wholesale via Range.InsertXML with a
# hand-built OOXML fragment (wrapped in the pkg:package form, which is what
# this host's InsertXML implementation expects).

$d = $word.ActiveDocument

function New-RunXml($segments) {
    # $segments is an array of "text|preserve|breaksAfter" strings.
    $inner = ""
    foreach ($seg in $segments) {
        $parts = $seg.Split("|")
        $text = $parts[0]
        $preserve = $parts[1]
        $breaks = [int]$parts[2]

        if ($preserve -eq "1") {
            $inner += '<w:t xml:space="preserve">' + $text + '</w:t>'
        } else {
            $inner += '<w:t>' + $text + '</w:t>'
        }
        for ($i = 0; $i -lt $breaks; $i++) {
            $inner += '<w:br/>'
        }
    }
    return $inner
}

function Find-ParagraphIndexByText($fullText) {
    # Paragraph.Range.Text carries a trailing paragraph-mark (CR); compare
    # against that so we need an exact, unambiguous match.
    $wanted = $fullText + "`r"
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        if ($t -eq $wanted) {
            return $i
        }
    }
    return -1
}

function Set-ParagraphRuns($originalText, $segments) {
    $paraIndex = Find-ParagraphIndexByText $originalText
    if ($paraIndex -eq -1) {
        throw "Could not locate paragraph with text: " + $originalText
    }
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $runInner = New-RunXml $segments

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData>' `
        + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        + '<w:body><w:p><w:r>' + $runInner + '</w:r></w:p></w:body>' `
        + '</w:document>' `
        + '</pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xml)
}

# --- "Programa resumido" paragraph ---------------------------------------
$orig1 = 'Fundamentos da Eletrônica; Sensores e atuadores; Microcontroladores;Sistemas mecatrônicos;Experimentos práticos.'
Set-ParagraphRuns $orig1 @(
    'Fundamentos da Eletrônica; |1|1',
    'Sensores e atuadores; |1|1',
    'Microcontroladores;|0|1',
    'Sistemas mecatrônicos;|0|1',
    'Experimentos práticos.|0|0'
)

# --- "Programa" paragraph --------------------------------------------------
$orig2 = 'Fundamentos da Eletrônica: digital e analógica; Uso de sensores e atuadores em processos produtivos; Programação básica de microcontroladores, com foco em Arduino ou similar;Projetos de uso em sistemas mecatrônicos aplicados a engenharia de produção;Experimentos práticos.'
Set-ParagraphRuns $orig2 @(
    'Fundamentos da Eletrônica: digital e analógica; |1|1',
    'Uso de sensores e atuadores em processos produtivos; |1|1',
    'Programação básica de microcontroladores, com foco em Arduino ou similar;|0|1',
    'Projetos de uso em sistemas mecatrônicos aplicados a engenharia de produção;|0|1',
    'Experimentos práticos.|0|0'
)

# --- "Bibliografia" paragraph ----------------------------------------------
$orig3 = 'Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U.B. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p.Bibliografia complementar será indicada ao longo do curso.'
Set-ParagraphRuns $orig3 @(
    'Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) |1|2',
    'BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. |1|2',
    'THOMAZINI, Daniel; ALBUQUERQUE, Pedro U.B. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p.|0|2',
    'Bibliografia complementar será indicada ao longo do curso.|0|0'
)

Write-Host "Done."
